# Update countries & provincias Spain
# Applies the sharedStrings reorder (reflected as cell-content swaps) and
# the refreshed case counts for the "Pais" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

$ws.Range('A1').Value = 'Datos actualizados a 28 de Junio de 2020 a las 16:37'
$ws.Range('B4').Value = 2597891
$ws.Range('C4').Value = 1354
$ws.Range('D4').Value = 1081531
$ws.Range('E4').Value = 1388197
$ws.Range('G4').Value = 11
$ws.Range('H4').Value = 128163
$ws.Range('B7').Value = 544615
$ws.Range('C7').Value = 15038
$ws.Range('D7').Value = 319393
$ws.Range('E7').Value = 208798
$ws.Range('G7').Value = 321
$ws.Range('H7').Value = 16424
$ws.Range('B8').Value = 311151
$ws.Range('C8').Value = 901
$ws.Range('G8').Value = 36
$ws.Range('H8').Value = 43550
$ws.Range('B18').Value = 182493
$ws.Range('C18').Value = 3989
$ws.Range('E18').Value = 56187
$ws.Range('D30').Value = 20134
$ws.Range('E30').Value = 36393
$ws.Range('G30').Value = 10
$ws.Range('H30').Value = 1217
$ws.Range('A35').Value = 'Irak'
$ws.Range('B35').Value = 45402
$ws.Range('C35').Value = 2140
$ws.Range('D35').Value = 21122
$ws.Range('E35').Value = 22524
$ws.Range('G35').Value = 96
$ws.Range('H35').Value = 1756
$ws.Range('A36').Value = 'Kuwait'
$ws.Range('B36').Value = 44942
$ws.Range('C36').Value = 551
$ws.Range('D36').Value = 35494
$ws.Range('E36').Value = 9100
$ws.Range('G36').Value = 4
$ws.Range('H36').Value = 348
$ws.Range('A37').Value = 'Singapur'
$ws.Range('B37').Value = 43459
$ws.Range('C37').Value = 213
$ws.Range('D37').Value = 37163
$ws.Range('E37').Value = 6270
$ws.Range('H37').Value = 26
$ws.Range('B60').Value = 16250
$ws.Range('C60').Value = 170
$ws.Range('E60').Value = 6639
$ws.Range('G60').Value = 9
$ws.Range('H60').Value = 530
$ws.Range('B62').Value = 14046
$ws.Range('C62').Value = 254
$ws.Range('D62').Value = 12464
$ws.Range('E62').Value = 1312
$ws.Range('G62').Value = 3
$ws.Range('H62').Value = 270
$ws.Range('E74').Value = 2544
$ws.Range('G74').Value = 1
$ws.Range('H74').Value = 21
$ws.Range('A79').Value = 'Republica de Macedonia'
$ws.Range('B79').Value = 6080
$ws.Range('C79').Value = 174
$ws.Range('D79').Value = 2315
$ws.Range('E79').Value = 3479
$ws.Range('H79').Value = 286
$ws.Range('A80').Value = 'Kenia'
$ws.Range('B80').Value = 6070
$ws.Range('C80').Value = 259
$ws.Range('D80').Value = 1936
$ws.Range('E80').Value = 3993
$ws.Range('H80').Value = 141
$ws.Range('A81').Value = 'El Salvador'
$ws.Range('B81').Value = 5934
$ws.Range('C81').Value = 207
$ws.Range('D81').Value = 3557
$ws.Range('E81').Value = 2225
$ws.Range('G81').Value = 9
$ws.Range('H81').Value = 152
$ws.Range('B82').Value = 5849
$ws.Range('C82').Value = 50
$ws.Range('D82').Value = 4448
$ws.Range('E82').Value = 1349
$ws.Range('B84').Value = 5689
$ws.Range('C84').Value = 119
$ws.Range('D84').Value = 2132
$ws.Range('E84').Value = 3459
$ws.Range('G84').Value = 4
$ws.Range('H84').Value = 98
$ws.Range('B104').Value = 2332
$ws.Range('C104').Value = 2
$ws.Range('D104').Value = 2201
$ws.Range('E104').Value = 45
$ws.Range('E120').Value = 86
$ws.Range('G120').Value = 2
$ws.Range('H120').Value = 111
$ws.Range('A205').Value = 'Dominica'
$ws.Range('A206').Value = 'Fiyi'
$ws.Range('A209').Value = 'Islas Malvinas'
$ws.Range('A210').Value = 'Groenlandia'
